$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 15999.667
$ws.Range("J70").Value = 15999.667
$ws.Range("L70").Value = 47999.001
$ws.Range("N70").Value = -48539.001
$ws.Range("H73").Value = 15999.667
$ws.Range("J73").Value = 15999.667
$ws.Range("L73").Value = 47999.001
$ws.Range("N73").Value = -49871.001
$ws.Range("H132").Value = 1035.5807
$ws.Range("I132").Value = 1035.5807
$ws.Range("K132").Value = 3106.7421
$ws.Range("M132").Value = -576.7420999999999
$ws.Range("H138").Value = 1899.2222
$ws.Range("I138").Value = 1308.091
$ws.Range("J138").Value = 4500.2
$ws.Range("K138").Value = 3924.273
$ws.Range("L138").Value = 13500.6
$ws.Range("M138").Value = 1215.727
$ws.Range("N138").Value = -23780.6

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4273
$ws.Range("I63").Value = 3114
$ws.Range("J63").Value = 7750
$ws.Range("K63").Value = 3114
$ws.Range("L63").Value = 7750
$ws.Range("M63").Value = -2428
$ws.Range("N63").Value = -9122
$ws.Range("H66").Value = 4273
$ws.Range("I66").Value = 3114
$ws.Range("J66").Value = 7750
$ws.Range("K66").Value = 15570
$ws.Range("L66").Value = 38750
$ws.Range("M66").Value = -12138
$ws.Range("N66").Value = -45614
$ws.Range("H104").Value = 34999.5
$ws.Range("J104").Value = 34999.5
$ws.Range("L104").Value = 34999.5
$ws.Range("N104").Value = -41987.5
$ws.Range("H109").Value = 67198.78
$ws.Range("J109").Value = 67198.78
$ws.Range("L109").Value = 67198.78
$ws.Range("N109").Value = -69972.78
$ws.Range("H112").Value = 27997.4
$ws.Range("J112").Value = 27997.4
$ws.Range("L112").Value = 27997.4
$ws.Range("N112").Value = -30951.4
$ws.Range("H124").Value = 35900
$ws.Range("J124").Value = 35900
$ws.Range("L124").Value = 35900
$ws.Range("N124").Value = -45720
$ws.Range("H135").Value = 26555.6
$ws.Range("J135").Value = 26555.6
$ws.Range("L135").Value = 26555.6
$ws.Range("N135").Value = -36695.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 34282.43
$ws.Range("J81").Value = 34282.43
$ws.Range("L81").Value = 34282.43
$ws.Range("N81").Value = -36404.43
$ws.Range("H84").Value = 34282.43
$ws.Range("J84").Value = 34282.43
$ws.Range("L84").Value = 102847.29
$ws.Range("N84").Value = -113455.29
$ws.Range("H105").Value = 1871.8148
$ws.Range("I105").Value = 1861.091
$ws.Range("K105").Value = 1861.091
$ws.Range("M105").Value = -114.0909999999999
$ws.Range("H106").Value = 19000
$ws.Range("J106").Value = 19000
$ws.Range("L106").Value = 19000
$ws.Range("N106").Value = -21524
$ws.Range("H110").Value = 99992.75
$ws.Range("J110").Value = 99992.75
$ws.Range("L110").Value = 99992.75
$ws.Range("N110").Value = -108172.75
$ws.Range("H130").Value = 59994.2
$ws.Range("J130").Value = 59994.2
$ws.Range("L130").Value = 59994.2
$ws.Range("N130").Value = -70034.2
$ws.Range("H135").Value = 28784.857
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 28784.857
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 28784.857
$ws.Range("M135").ClearContents()
$ws.Range("N135").Value = -38924.857

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5437.6665
$ws.Range("J31").Value = 8188.9165
$ws.Range("L31").Value = 8188.9165
$ws.Range("N31").Value = -8778.916499999999
$ws.Range("H34").Value = 5437.6665
$ws.Range("J34").Value = 8188.9165
$ws.Range("L34").Value = 8188.9165
$ws.Range("N34").Value = -8592.916499999999
$ws.Range("H43").Value = 9499.5
$ws.Range("J43").Value = 9499.5
$ws.Range("L43").Value = 9499.5
$ws.Range("N43").Value = -9867.5
$ws.Range("H99").Value = 2237.25
$ws.Range("I99").Value = 2079.6
$ws.Range("K99").Value = 2079.6
$ws.Range("M99").Value = -581.5999999999999
$ws.Range("H101").Value = 9499.5
$ws.Range("J101").Value = 9499.5
$ws.Range("L101").Value = 9499.5
$ws.Range("N101").Value = -15989.5
$ws.Range("H122").Value = 2372.4285
$ws.Range("J122").Value = 2134.6667
$ws.Range("L122").Value = 6404.000100000001
$ws.Range("N122").Value = -11304.0001
$ws.Range("H126").Value = 2237.25
$ws.Range("I126").Value = 2079.6
$ws.Range("K126").Value = 6238.799999999999
$ws.Range("M126").Value = -3768.799999999999
$ws.Range("H132").Value = 2177.074
$ws.Range("I132").Value = 1060.6875
$ws.Range("K132").Value = 3182.0625
$ws.Range("M132").Value = -652.0625
$ws.Range("H134").Value = 2687.8572
$ws.Range("I134").Value = 2052.5833
$ws.Range("K134").Value = 6157.749899999999
$ws.Range("M134").Value = -3622.749899999999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H113").Value = 6593.2354
$ws.Range("I113").Value = 25450.25
$ws.Range("J113").Value = 791.0769
$ws.Range("K113").Value = 76350.75
$ws.Range("L113").Value = 2373.2307
$ws.Range("M113").Value = -74180.75
$ws.Range("N113").Value = -6713.2307
$ws.Range("H122").Value = 1278.8182
$ws.Range("J122").Value = 1677
$ws.Range("L122").Value = 15093
$ws.Range("N122").Value = -19993
$ws.Range("H134").Value = 6889.35
$ws.Range("I134").Value = 7210.6113
$ws.Range("K134").Value = 21631.8339
$ws.Range("M134").Value = -16561.8339
$ws.Range("H136").Value = 166668980
$ws.Range("I136").Value = 166668980
$ws.Range("K136").Value = 500006940
$ws.Range("M136").Value = -500001840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H98").Value = 9999.666999999999
$ws.Range("J98").Value = 9999.666999999999
$ws.Range("L98").Value = 9999.666999999999
$ws.Range("N98").Value = -15989.667
$ws.Range("H126").Value = 35015.324
$ws.Range("I126").Value = 2946.4092
$ws.Range("K126").Value = 8839.2276
$ws.Range("M126").Value = -6369.2276

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 31298.5
$ws.Range("I16").Value = 31298.5
$ws.Range("K16").Value = 31298.5
$ws.Range("M16").Value = -31128.5
$ws.Range("H104").Value = 9997
$ws.Range("J104").Value = 9997
$ws.Range("L104").Value = 9997
$ws.Range("N104").Value = -16985
$ws.Range("H110").Value = 22101.834
$ws.Range("J110").Value = 22101.834
$ws.Range("L110").Value = 22101.834
$ws.Range("N110").Value = -30281.834
$ws.Range("H127").Value = 49358
$ws.Range("J127").Value = 49358
$ws.Range("L127").Value = 49358
$ws.Range("N127").Value = -59278

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 56908.93
$ws.Range("I122").Value = 79180.2
$ws.Range("J122").Value = 1230.75
$ws.Range("K122").Value = 237540.6
$ws.Range("L122").Value = 3692.25
$ws.Range("M122").Value = -235090.6
$ws.Range("N122").Value = -8592.25
$ws.Range("H126").Value = 7834.069
$ws.Range("I126").Value = 8723.333000000001
$ws.Range("J126").Value = 5499.75
$ws.Range("K126").Value = 26169.999
$ws.Range("L126").Value = 16499.25
$ws.Range("M126").Value = -21439.25
